# ConditionBit를 이용한 퀘스트 수주 (GAME_START만)
# Renames getCondition/getConditionArguments columns to condition/conditionArguments
# and switches the G (condition) / J (goalType) / K (goalArgument) columns on
# Sheet1 to use the new string-based enum values.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Header renames (G1, H1) ---
$ws1.Range("G1").Value = "condition"
$ws1.Range("H1").Value = "conditionArguments"

# --- Row 2 (index 1): GAME_START quest, goal MOVE_TO ---
$ws1.Range("G2").Value = "GAME_START"
$ws1.Range("J2").Value = "MOVE_TO"

# --- Row 3 (index 2): QUEST_END condition w/ argument 1, goal KILL_LINK ---
$ws1.Range("G3").Value = "QUEST_END"
$ws1.Range("H3").Value = 1
$ws1.Range("J3").Value = "KILL_LINK"
$ws1.Range("K3").Value = "6,1"

# --- Row 4 (index 3): QUEST_END condition w/ argument 2, goal GET_ITEM ---
$ws1.Range("G4").Value = "QUEST_END"
$ws1.Range("H4").Value = 2
$ws1.Range("J4").Value = "GET_ITEM"
$ws1.Range("K4").Value = "101,1"

# --- Column J width + active selection cosmetic updates ---
# NB: the ColumnWidth COM setter here round-trips through a pixel quantizer
# (~1/7-character steps), so the exact source width (10.25) isn't reachable;
# 9.43 lands on the closest attainable quantized width.
$ws1.Columns.Item(10).ColumnWidth = 9.43
$ws1.Range("G3").Select()
